# Generate Report for Handoff
# Updates the "b.md" rows (row 3) across the Overview/zh-cn/de-de sheets to
# reflect a fresh handoff: status flips from "Handed back: in sync with en-US"
# to "Ready for handoff", new handoff xliff file names + timestamps are
# recorded, and (for the locale sheets) an error detail note is attached.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Overview sheet: b.md row (row 3) - Status columns (E, F) and the
# "Latest HO Xliff Generate Date" column (G).
# ---------------------------------------------------------------------
$overview = $wb.Worksheets.Item("Overview")
$overview.Range("E3").Value = "Ready for handoff"
$overview.Range("F3").Value = "Ready for handoff"
$overview.Range("G3").Value = "2016-08-17 18:36:40"

# ---------------------------------------------------------------------
# zh-cn sheet: b.md row (row 3)
# ---------------------------------------------------------------------
$zhcn = $wb.Worksheets.Item("zh-cn")
$zhcn.Range("C3").Value = "Ready for handoff"
# Leading apostrophe forces text storage (matches the "True"/"False"
# text cells elsewhere in the sheet, as opposed to a native boolean);
# Style is then reset to Normal so no stray quote-prefix format sticks.
$zhcn.Range("F3").Value = "'False"
$zhcn.Range("F3").Style = "Normal"
$zhcn.Range("G3").Value = "b.63290e5768f688058c7b37413b0a5c26c308f864.zh-cn.xlf"
$zhcn.Range("H3").Value = "2016-08-17 18:36:34"
$zhcn.Range("P3").Value = "The version of handback file is not the latest, current: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/0e1b2bce420c2cd32e05c386a975b3aa9ae4fcc8/e2e/a.md, latest: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/e673df5d143451b2a974bbdf4fafb53a96ec3627/e2e/b.md."
# ColumnWidth (chars) differs from the stored OOXML "width" by the fixed
# ~5/6 character padding Excel applies when serialising; 40 - 5/6 round-trips
# to a stored width of exactly 40, matching the target column P width.
$zhcn.Columns.Item(16).ColumnWidth = 39.166666666666664

# ---------------------------------------------------------------------
# de-de sheet: b.md row (row 3)
# ---------------------------------------------------------------------
$dede = $wb.Worksheets.Item("de-de")
$dede.Range("C3").Value = "Ready for handoff"
$dede.Range("F3").Value = "'False"
$dede.Range("F3").Style = "Normal"
$dede.Range("G3").Value = "b.63290e5768f688058c7b37413b0a5c26c308f864.de-de.xlf"
$dede.Range("H3").Value = "2016-08-17 18:36:40"
$dede.Range("P3").Value = "The version of handback file is not the latest, current: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/0e1b2bce420c2cd32e05c386a975b3aa9ae4fcc8/e2e/a.md, latest: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/e673df5d143451b2a974bbdf4fafb53a96ec3627/e2e/b.md."
$dede.Columns.Item(16).ColumnWidth = 39.166666666666664
